$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 66, shifting existing rows 66-106 down to 67-107.
$ws.Rows.Item(66).Insert()

# Populate the newly inserted row 66 with the new weekly data point.
$ws.Cells.Item(66, 1).Value2 = 10
$ws.Cells.Item(66, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(66, 3).Value2 = "La Araucanía"
$ws.Cells.Item(66, 4).Value2 = 45072
$ws.Cells.Item(66, 5).Value2 = 9
$ws.Cells.Item(66, 6).Value2 = 300000001
$ws.Cells.Item(66, 7).Value2 = "Rabanito"
$ws.Cells.Item(66, 8).Value2 = "Sin especificar"
$ws.Cells.Item(66, 9).Value2 = "Primera"
$ws.Cells.Item(66, 10).Value2 = 50
$ws.Cells.Item(66, 11).Value2 = 8000
$ws.Cells.Item(66, 12).Value2 = 8000
$ws.Cells.Item(66, 13).Value2 = 8000
$ws.Cells.Item(66, 14).Value2 = "`$/docena de paquetes"
$ws.Cells.Item(66, 15).Value2 = "Provincia de Cautín"
$ws.Cells.Item(66, 16).Value2 = 667
$ws.Cells.Item(66, 17).Value2 = 12
$ws.Cells.Item(66, 18).Value2 = "Hortaliza"
